$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: C12 changes from a blank-ish text marker to a real CMP-updated-at date ---
$ws.Range("C12").Value = 46058.4680555556

# --- Row 22: first "closed" row gets promoted to a real UUID entry with a distinct
#     (Cascadia Mono SemiBold) font + wrap, and a " vbms" value in column B ---
$ws.Range("A22").Value = "c7010f85-2a8b-47e2-bd99-3a45bff61b4e"
$ws.Range("A22").Font.Name = "Cascadia Mono SemiBold"
$ws.Range("A22").WrapText = $true
$ws.Range("A22").HorizontalAlignment = 1
$ws.Range("B22").Value = " vbms"

# --- Row 23: gets its own UUID + " vbms" (style already matches target) ---
$ws.Range("A23").Value = "fa5b0a1f-2ca3-4ed2-9e8f-21d29b188ad8"
$ws.Range("B23").Value = " vbms"

# --- Rows 24-29: previously-empty rows now populated with UUID / vbms / date ---
$ws.Range("A24").Value = "056cc35c-d270-4f62-b696-65fced00fd51"
$ws.Range("B24").Value = " vbms"
$ws.Range("C24").Value = 46064.7840277778

$ws.Range("A25").Value = "87efed84-a179-49ec-9d0d-c4d29c9b2a1d"
$ws.Range("B25").Value = " vbms"
$ws.Range("C25").Value = 46065.4166666667

$ws.Range("A26").Value = "b8c73c9c-1edb-4a7a-aff1-c20ea37e5816"
$ws.Range("B26").Value = " vbms"
$ws.Range("C26").Value = 46065.5527777778

$ws.Range("A27").Value = "3e874961-c725-40be-861f-09b1fc3511d6"
$ws.Range("B27").Value = " vbms"
$ws.Range("C27").Value = 46065.725

$ws.Range("A28").Value = "26c58d71-8daf-460d-b6eb-59b4f4974912"
$ws.Range("B28").Value = " vbms"
$ws.Range("C28").Value = 46034.8090277778

$ws.Range("A29").Value = "224f98e7-98b6-41ed-9366-5f39d1a3409f"
$ws.Range("B29").Value = " vbms"
$ws.Range("C29").Value = 46066.3798611111

# --- Recalc so the dependent summary formulas (B4/C4) pick up the newly
#     populated rows, then move the selection to the last cell touched ---
$excel.Calculate()
$ws.Range("C29").Select()
